$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Personnes")

$ws.Cells.Item(4, 1).Value = 3.0
$ws.Cells.Item(4, 2).Value = "bernardghgh"
$ws.Cells.Item(4, 3).NumberFormat = "General"

$ws.Cells.Item(5, 1).Value = 4.0
$ws.Cells.Item(5, 2).Value = "bernardfhdkfh"
$ws.Cells.Item(5, 3).NumberFormat = "General"
